# Subscript industry CCS by energy related vs. process emissions (#129)
#
# The "CPPbI" (CPP CO2 Capture Potential by Industry) sheet previously had a
# single data column (B) giving the capturable share per industry. This
# change relabels that as "energy related emissions" and adds a second,
# identical data column C labelled "process emissions".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPPbI")

# New column headers (row 1) for the two emission-type sub-columns.
$ws.Range("B1").Value = "energy related emissions"
$ws.Range("C1").Value = "process emissions"

# Mirror column B's values into the new column C for every industry row.
for ($row = 2; $row -le 9; $row++) {
    $b = $ws.Cells.Item($row, 2)
    $c = $ws.Cells.Item($row, 3)
    $c.Value2 = $b.Value2
}

# Widen column A to fit the longer industry labels, and size the two new
# data columns.
$ws.Columns.Item(1).ColumnWidth = 42.833333333333336
$ws.Columns.Item(2).ColumnWidth = 23.833333333333336
$ws.Columns.Item(3).ColumnWidth = 25.0
